# Update "想去人数" (column F) values on 展览 / 演出 / 全部类型 sheets
# per gh-pages regeneration commit 456a3b4.
$wb = $excel.ActiveWorkbook

# --- 展览 (sheet 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(5, 6).Value = 1047
$ws.Cells.Item(7, 6).Value = 2695
$ws.Cells.Item(8, 6).Value = 42
$ws.Cells.Item(9, 6).Value = 1323
$ws.Cells.Item(10, 6).Value = 943
$ws.Cells.Item(11, 6).Value = 636
$ws.Cells.Item(13, 6).Value = 1203
$ws.Cells.Item(16, 6).Value = 756
$ws.Cells.Item(19, 6).Value = 543
$ws.Cells.Item(22, 6).Value = 660
$ws.Cells.Item(23, 6).Value = 618
$ws.Cells.Item(28, 6).Value = 625
$ws.Cells.Item(29, 6).Value = 6842
$ws.Cells.Item(34, 6).Value = 190
$ws.Cells.Item(35, 6).Value = 1657
$ws.Cells.Item(37, 6).Value = 115
$ws.Cells.Item(39, 6).Value = 148
$ws.Cells.Item(41, 6).Value = 156
$ws.Cells.Item(42, 6).Value = 21
$ws.Cells.Item(43, 6).Value = 79
$ws.Cells.Item(45, 6).Value = 150
$ws.Cells.Item(46, 6).Value = 141
$ws.Cells.Item(47, 6).Value = 125

# --- 演出 (sheet 2) ---
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(7, 6).Value = 3
$ws.Cells.Item(9, 6).Value = 56
$ws.Cells.Item(12, 6).Value = 202
$ws.Cells.Item(13, 6).Value = 4416
$ws.Cells.Item(14, 6).Value = 43
$ws.Cells.Item(18, 6).Value = 220

# --- 全部类型 (sheet 4) ---
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(6, 6).Value = 1047
$ws.Cells.Item(7, 6).Value = 2695
$ws.Cells.Item(8, 6).Value = 42
$ws.Cells.Item(9, 6).Value = 1323
$ws.Cells.Item(10, 6).Value = 943
$ws.Cells.Item(11, 6).Value = 636
$ws.Cells.Item(13, 6).Value = 1203
$ws.Cells.Item(17, 6).Value = 756
$ws.Cells.Item(21, 6).Value = 543
$ws.Cells.Item(23, 6).Value = 3
$ws.Cells.Item(24, 6).Value = 56
$ws.Cells.Item(25, 6).Value = 660
$ws.Cells.Item(26, 6).Value = 618
$ws.Cells.Item(30, 6).Value = 625
$ws.Cells.Item(31, 6).Value = 6842
$ws.Cells.Item(32, 6).Value = 202
$ws.Cells.Item(36, 6).Value = 190
$ws.Cells.Item(37, 6).Value = 1657
$ws.Cells.Item(40, 6).Value = 43
$ws.Cells.Item(41, 6).Value = 43
$ws.Cells.Item(44, 6).Value = 21
$ws.Cells.Item(45, 6).Value = 79
$ws.Cells.Item(46, 6).Value = 150
$ws.Cells.Item(48, 6).Value = 125
